# Append two new days of quote data (rows 64 and 65) to Sheet1, matching
# the style/number-format already used by the existing date column (A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (incl. the date number format on column A) from the
# last existing data row (63) down into the two new rows, so the new cells
# reuse the same style instead of creating new ones.
$ws.Range("A63:D63").Copy($ws.Range("A64:D64"))
$ws.Range("A63:D63").Copy($ws.Range("A65:D65"))

# Row 64 -> 2022-08-01 (serial 44774)
$ws.Range("A64").Value = 44774
$ws.Range("B64").Value = 410.7699890136719
$ws.Range("C64").Value = 67.41000366210938
$ws.Range("D64").Value = 74.02999877929688

# Row 65 -> 2022-08-02 (serial 44775)
$ws.Range("A65").Value = 44775
$ws.Range("B65").Value = 411.7550048828125
$ws.Range("C65").Value = 67.52999877929688
$ws.Range("D65").Value = 74.47000122070312
